# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45181 to 45182 (i.e. bump the date forward by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 221
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
